$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -87.8714
$ws.Range("B2").Value = -87.8618

$ws.Range("A3").Value = 30.2928
$ws.Range("B3").Value = 30.301

$ws.Range("A4").Value = -87.7507
$ws.Range("B4").Value = -87.7603

$ws.Range("A5").Value = 30.4921
$ws.Range("B5").Value = 30.4838
